# Update Hgf-Sdc2 NATMI edge-weight table with recomputed TPM-based statistics.
# (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06674466666666666
$ws.Range("H2").Value = 0.200234
$ws.Range("I2").Value = 0.0009912440954723497
$ws.Range("J2").Value = 0.0009958565080158308
$ws.Range("M2").Value = 1.009860666666667
$ws.Range("N2").Value = 3.029582
$ws.Range("O2").Value = 0.01353413605720072
$ws.Range("P2").Value = 0.01542521070970148
$ws.Range("Q2").Value = 0.06740281357644444
$ws.Range("R2").Value = 0.606625322188
$ws.Range("S2").Value = [double]"1.341563245401964e-05"
$ws.Range("T2").Value = [double]"1.536129647277171e-05"

# Row 3: ECs -> FAPs
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06674466666666666
$ws.Range("H3").Value = 0.200234
$ws.Range("I3").Value = 0.0009912440954723497
$ws.Range("J3").Value = 0.0009958565080158308
$ws.Range("O3").Value = 0.6185519418990597
$ws.Range("P3").Value = 0.704979911415303
$ws.Range("Q3").Value = 3.080517371109777
$ws.Range("R3").Value = 27.72465633998799
$ws.Range("S3").Value = 0.0006131359601503989
$ws.Range("T3").Value = 0.0007020588328033534

# Row 4: ECs -> Inflammatory-Mac
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06674466666666666
$ws.Range("H4").Value = 0.200234
$ws.Range("I4").Value = 0.0009912440954723497
$ws.Range("J4").Value = 0.0009958565080158308
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.009315666666666667
$ws.Range("N4").Value = 0.027947
$ws.Range("O4").Value = 0.0001248484115599408
$ws.Range("P4").Value = 0.000142293017222847
$ws.Range("Q4").Value = 0.0006217710664444444
$ws.Range("R4").Value = 0.005595939598
$ws.Range("S4").Value = [double]"1.237552507878932e-07"
$ws.Range("T4").Value = [double]"1.417034272465809e-07"

# Row 5: ECs -> MuSCs
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.06674466666666666
$ws.Range("H5").Value = 0.200234
$ws.Range("I5").Value = 0.0009912440954723497
$ws.Range("J5").Value = 0.0009958565080158308
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 27.4428835
$ws.Range("N5").Value = 54.885767
$ws.Range("O5").Value = 0.3677890736321797
$ws.Range("P5").Value = 0.2794525848577725
$ws.Range("Q5").Value = 1.831666111579667
$ws.Range("R5").Value = 10.989996669478
$ws.Range("S5").Value = 0.0003645687476171433
$ws.Range("T5").Value = 0.000278294675312459

# Row 6: FAPs -> ECs
$ws.Range("I6").Value = 0.1187608236941705
$ws.Range("J6").Value = 0.1193134362296531
$ws.Range("M6").Value = 1.009860666666667
$ws.Range("N6").Value = 3.029582
$ws.Range("O6").Value = 0.01353413605720072
$ws.Range("P6").Value = 0.01542521070970148
$ws.Range("Q6").Value = 8.075522160693112
$ws.Range("R6").Value = 72.679699446238
$ws.Range("S6").Value = 0.00160732514614213
$ws.Range("T6").Value = 0.001840434894340929

# Row 7: FAPs -> FAPs
$ws.Range("I7").Value = 0.1187608236941705
$ws.Range("J7").Value = 0.1193134362296531
$ws.Range("O7").Value = 0.6185519418990597
$ws.Range("P7").Value = 0.704979911415303
$ws.Range("S7").Value = 0.07345973811756101
$ws.Range("T7").Value = 0.08411357570383622

# Row 8: FAPs -> Inflammatory-Mac
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.1187608236941705
$ws.Range("J8").Value = 0.1193134362296531
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.009315666666666667
$ws.Range("N8").Value = 0.027947
$ws.Range("O8").Value = 0.0001248484115599408
$ws.Range("P8").Value = 0.000142293017222847
$ws.Range("Q8").Value = 0.07449430905811111
$ws.Range("R8").Value = 0.670448781523
$ws.Range("S8").Value = [double]"1.482710019376736e-05"
$ws.Range("T8").Value = [double]"1.697746883634308e-05"

# Row 9: FAPs -> MuSCs
$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.1187608236941705
$ws.Range("J9").Value = 0.1193134362296531
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 27.4428835
$ws.Range("N9").Value = 54.885767
$ws.Range("O9").Value = 0.3677890736321797
$ws.Range("P9").Value = 0.2794525848577725
$ws.Range("Q9").Value = 219.4516740503172
$ws.Range("R9").Value = 1316.710044301903
$ws.Range("S9").Value = 0.04367893333027357
$ws.Range("T9").Value = 0.03334244816263956

# Row 10: Inflammatory-Mac -> ECs
$ws.Range("G10").Value = 24.06383433333333
$ws.Range("H10").Value = 72.191503
$ws.Range("I10").Value = 0.3573788721796719
$ws.Range("J10").Value = 0.3590418115105046
$ws.Range("M10").Value = 1.009860666666667
$ws.Range("N10").Value = 3.029582
$ws.Range("O10").Value = 0.01353413605720072
$ws.Range("P10").Value = 0.01542521070970148
$ws.Range("Q10").Value = 24.30111978241622
$ws.Range("R10").Value = 218.710078041746
$ws.Range("S10").Value = 0.004836814280048626
$ws.Range("T10").Value = 0.005538295596142455

# Row 11: Inflammatory-Mac -> FAPs
$ws.Range("G11").Value = 24.06383433333333
$ws.Range("H11").Value = 72.191503
$ws.Range("I11").Value = 0.3573788721796719
$ws.Range("J11").Value = 0.3590418115105046
$ws.Range("O11").Value = 0.6185519418990597
$ws.Range("P11").Value = 0.704979911415303
$ws.Range("Q11").Value = 1110.636450542983
$ws.Range("R11").Value = 9995.728054886844
$ws.Range("S11").Value = 0.2210573953804319
$ws.Range("T11").Value = 0.2531172644730654

# Row 12: Inflammatory-Mac -> Inflammatory-Mac
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 24.06383433333333
$ws.Range("H12").Value = 72.191503
$ws.Range("I12").Value = 0.3573788721796719
$ws.Range("J12").Value = 0.3590418115105046
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.009315666666666667
$ws.Range("N12").Value = 0.027947
$ws.Range("O12").Value = 0.0001248484115599408
$ws.Range("P12").Value = 0.000142293017222847
$ws.Range("Q12").Value = 0.2241706593712222
$ws.Range("R12").Value = 2.017535934341
$ws.Range("S12").Value = [double]"4.461818451671515e-05"
$ws.Range("T12").Value = [double]"5.108914266898641e-05"

# Row 13: Inflammatory-Mac -> MuSCs
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 24.06383433333333
$ws.Range("H13").Value = 72.191503
$ws.Range("I13").Value = 0.3573788721796719
$ws.Range("J13").Value = 0.3590418115105046
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 27.4428835
$ws.Range("N13").Value = 54.885767
$ws.Range("O13").Value = 0.3677890736321797
$ws.Range("P13").Value = 0.2794525848577725
$ws.Range("Q13").Value = 660.3810021729669
$ws.Range("R13").Value = 3962.286013037801
$ws.Range("S13").Value = 0.1314400443346747
$ws.Range("T13").Value = 0.1003351622986276

# Row 14: MuSCs -> ECs
$ws.Range("G14").Value = 0.9355965
$ws.Range("H14").Value = 1.871193
$ws.Range("I14").Value = 0.01389481066706348
$ws.Range("J14").Value = 0.009306310251024633
$ws.Range("M14").Value = 1.009860666666667
$ws.Range("N14").Value = 3.029582
$ws.Range("O14").Value = 0.01353413605720072
$ws.Range("P14").Value = 0.01542521070970148
$ws.Range("Q14").Value = 0.9448221052210001
$ws.Range("R14").Value = 5.668932631326
$ws.Range("S14").Value = 0.0001880542580570811
$ws.Range("T14").Value = 0.0001435517965519098

# Row 15: MuSCs -> FAPs
$ws.Range("G15").Value = 0.9355965
$ws.Range("H15").Value = 1.871193
$ws.Range("I15").Value = 0.01389481066706348
$ws.Range("J15").Value = 0.009306310251024633
$ws.Range("O15").Value = 0.6185519418990597
$ws.Range("P15").Value = 0.704979911415303
$ws.Range("Q15").Value = 43.181296941571
$ws.Range("R15").Value = 259.087781649426
$ws.Range("S15").Value = 0.008594662120431886
$ws.Range("T15").Value = 0.006560761776370673

# Row 16: MuSCs -> Inflammatory-Mac
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = 0.9355965
$ws.Range("H16").Value = 1.871193
$ws.Range("I16").Value = 0.01389481066706348
$ws.Range("J16").Value = 0.009306310251024633
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.009315666666666667
$ws.Range("N16").Value = 0.027947
$ws.Range("O16").Value = 0.0001248484115599408
$ws.Range("P16").Value = 0.000142293017222847
$ws.Range("Q16").Value = 0.0087157051285
$ws.Range("R16").Value = 0.052294230771
$ws.Range("S16").Value = [double]"1.734745040708997e-06"
$ws.Range("T16").Value = [double]"1.324222964830206e-06"

# Row 17: MuSCs -> MuSCs
$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = 0.9355965
$ws.Range("H17").Value = 1.871193
$ws.Range("I17").Value = 0.01389481066706348
$ws.Range("J17").Value = 0.009306310251024633
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 27.4428835
$ws.Range("N17").Value = 54.885767
$ws.Range("O17").Value = 0.3677890736321797
$ws.Range("P17").Value = 0.2794525848577725
$ws.Range("Q17").Value = 25.67546575250775
$ws.Range("R17").Value = 102.701863010031
$ws.Range("S17").Value = 0.005110359543533807
$ws.Range("T17").Value = 0.00260067245513722

# Row 18: Resolving-Mac -> ECs
$ws.Range("G18").Value = 34.27139366666666
$ws.Range("H18").Value = 102.814181
$ws.Range("I18").Value = 0.5089742493636218
$ws.Range("J18").Value = 0.5113425855008019
$ws.Range("M18").Value = 1.009860666666667
$ws.Range("N18").Value = 3.029582
$ws.Range("O18").Value = 0.01353413605720072
$ws.Range("P18").Value = 0.01542521070970148
$ws.Range("Q18").Value = 34.60933245581577
$ws.Range("R18").Value = 311.483992102342
$ws.Range("S18").Value = 0.006888526740498865
$ws.Range("T18").Value = 0.007887567126193415

# Row 19: Resolving-Mac -> FAPs
$ws.Range("G19").Value = 34.27139366666666
$ws.Range("H19").Value = 102.814181
$ws.Range("I19").Value = 0.5089742493636218
$ws.Range("J19").Value = 0.5113425855008019
$ws.Range("O19").Value = 0.6185519418990597
$ws.Range("P19").Value = 0.704979911415303
$ws.Range("Q19").Value = 1581.753701004449
$ws.Range("R19").Value = 14235.78330904004
$ws.Range("S19").Value = 0.3148270103204845
$ws.Range("T19").Value = 0.3604862506292273

# Row 20: Resolving-Mac -> Inflammatory-Mac
$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("G20").Value = 34.27139366666666
$ws.Range("H20").Value = 102.814181
$ws.Range("I20").Value = 0.5089742493636218
$ws.Range("J20").Value = 0.5113425855008019
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.009315666666666667
$ws.Range("N20").Value = 0.027947
$ws.Range("O20").Value = 0.0001248484115599408
$ws.Range("P20").Value = 0.000142293017222847
$ws.Range("Q20").Value = 0.3192608796007777
$ws.Range("R20").Value = 2.873347916406999
$ws.Range("S20").Value = [double]"6.354462655796138e-05"
$ws.Range("T20").Value = [double]"7.276047932544072e-05"

# Row 21: Resolving-Mac -> MuSCs
$ws.Range("D21").Value = "MuSCs"
$ws.Range("G21").Value = 34.27139366666666
$ws.Range("H21").Value = 102.814181
$ws.Range("I21").Value = 0.5089742493636218
$ws.Range("J21").Value = 0.5113425855008019
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 27.4428835
$ws.Range("N21").Value = 54.885767
$ws.Range("O21").Value = 0.3677890736321797
$ws.Range("P21").Value = 0.2794525848577725
$ws.Range("Q21").Value = 940.5058637769711
$ws.Range("R21").Value = 5643.035182661826
$ws.Range("S21").Value = 0.1871951676760805
$ws.Range("T21").Value = 0.1428960072660556
